# Add the new "ManageNews" worksheet after the last existing sheet (AdminPage)
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ManageNews"

# Populate A1 with the news headline text (becomes a new shared string)
$newSheet.Range("A1").Value = "News Update - Glen"

# Match column A's width to the source layout (~23.42 chars)
$newSheet.Columns.Item(1).ColumnWidth = 22.6

# Mirror the page setup used on the other sheets (paper size, fit-to-page,
# portrait orientation, and the same margins as AdminPage)
$ps = $newSheet.PageSetup
$ps.PaperSize = 9
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.Orientation = 1
$ps.LeftMargin = 50.45669291338584
$ps.RightMargin = 50.45669291338584
$ps.TopMargin = 54.14173228346456
$ps.BottomMargin = 54.14173228346456
$ps.HeaderMargin = 21.599999999999998
$ps.FooterMargin = 21.599999999999998

# Make ManageNews the active tab (activeTab becomes 2, 0-based, third sheet)
$newSheet.Activate()
